# STATS.xlsx - "new comparaison with Jaro-Winkler algo"
#
# Updates the Surface comparison block (row 17, newly used as a note row
# like row 22) and refreshes several success/failure counts that feed the
# totals in rows 24/25/27, after a re-run of the comparison with a
# Jaro-Winkler based matcher.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (LIGNE C) : refreshed success/failure counts ---------------
$ws.Range("J6").Value = 7802
$ws.Range("K6").Value = 4

# --- Row 15 (T5) : refreshed success/failure counts --------------------
$ws.Range("J15").Value = 776
$ws.Range("K15").Value = 3

# --- Row 17 (LIAISONS B/D) : turned into a note row, like row 22 -------
# Copy the look (fill/merge style) of the existing note row 22 onto row 17
# before changing its content, so B17:G17 gets style 29 on every cell.
$ws.Range("B22:G22").Copy()
$ws.Range("B17:G17").PasteSpecial(-4122)
$ws.Range("B17:G17").Merge()

# Refresh the comment text (row 22's text also changes - see below) and
# put the new note in the now-styled B17.
$ws.Range("B22").Value = "Surface vs surface 3h43 pour 20498 succès et 42518 échecs,  après reprise des échecs sur tout hors ligne fortes :"
$ws.Range("B17").Value = "En comparant les échecs à tout darfeuille on passe de 14842 échecs à 11388"

$ws.Range("H17").Value = 100666
$ws.Range("J17").Formula = "=42780+3456"
$ws.Range("K17").Value = 11388

# --- Row 18 (TRAM COMMUNS) : refreshed success/failure counts ----------
$ws.Range("J18").Value = 6565
$ws.Range("K18").Value = 30

# --- Row 22 (Surface) : refreshed success/failure counts ---------------
$ws.Range("H22").Value = 106666
$ws.Range("J22").Value = 38256
$ws.Range("K22").Value = 24660

# --- Row 24 (TOTAL) : H24 now sums through row 23, not just row 19 -----
$ws.Range("H24").Formula = "=SUM(H4:H23)"

# --- Conditional formatting on column G : row 17 is now a note row, so
# exclude it from the percentage colour scale (it used to span G1:G21).
$gRuleRange = $ws.Range("G1:G21,G26,G28:G1048576,G23:G24")
$gRule = $gRuleRange.FormatConditions.Item(1)
$gRule.ModifyAppliesToRange($ws.Range("G1:G16"))

# --- Selection cursor, as left by the author ----------------------------
$ws.Range("I24").Select()

Write-Host "edit applied"
